$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new diary entries (rows 29 and 30)
# Fill column A first (dates), then column B (descriptions), to match
# the shared-string insertion order of the original edit.
$ws.Range("A29").Value = "24 марта"
$ws.Range("A30").Value = "25 марта"

$ws.Range("B30").Value = "Окончательный рефакторинг проекта бенчмаров и оформление проекта"
$ws.Range("B29").Value = "Написание скрипта на python для построения графика по данным результатов работы бенчмарка"

# Match the right-aligned format used throughout column A/B on this sheet
$ws.Range("B29:B30").HorizontalAlignment = -4152

# Update selection to match the new active cell
$ws.Range("B31").Select()
